$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header row (row 1) ---
$ws.Range("F1").Value = "Roles > Group Id"
$ws.Range("G1").Value = "Roles > User Id"

# --- Row 2: drop old Data Steward user name, add new group id value ---
$ws.Range("F2").Clear()
$ws.Range("G2").Value = "00000000-0000-0000-0000-000000900003"
$ws.Range("I2").Clear()

# --- Row 3: drop old Data Steward group id, add new group id value in F ---
$ws.Range("G3").Clear()
$ws.Range("F3").Value = "00000000-0000-0000-0000-000001000003"

# --- Row 4: drop the reviewer/group name value entirely ---
$ws.Range("H4").Clear()

# --- Remove now unused columns H and I (shrinks used range to A1:G4) ---
$ws.Range("H1:I4").Clear()
$ws.Columns("H:I").ColumnWidth = 8.83203125

# --- Re-apply explicit column widths for the two responsibility columns ---
$ws.Columns("F").ColumnWidth = 33.5
$ws.Columns("G").ColumnWidth = 35

# --- Selection matches the authored workbook ---
$ws.Range("F3").Select()
